$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Cases")
$ws.Cells.Item(5, 28).Value = 13
$ws.Cells.Item(5, 9).Value = 4
$ws.Cells.Item(6, 28).Value = 20
$ws.Cells.Item(6, 9).Value = 6
$ws.Cells.Item(7, 28).Value = 27
$ws.Cells.Item(7, 9).Value = 8
$ws.Cells.Item(8, 28).Value = 37
$ws.Cells.Item(8, 9).Value = 8
$ws.Cells.Item(9, 28).Value = 54
$ws.Cells.Item(9, 9).Value = 11
$ws.Cells.Item(10, 28).Value = 73
$ws.Cells.Item(10, 9).Value = 11
$ws.Cells.Item(11, 28).Value = 118
$ws.Cells.Item(11, 9).Value = 15
$ws.Cells.Item(12, 28).Value = 182
$ws.Cells.Item(12, 9).Value = 19
$ws.Cells.Item(14, 28).Value = 306
$ws.Cells.Item(14, 9).Value = 38
$ws.Cells.Item(15, 28).Value = 383
$ws.Cells.Item(15, 9).Value = 47
$ws.Cells.Item(16, 28).Value = 520
$ws.Cells.Item(16, 9).Value = 73
$ws.Cells.Item(17, 28).Value = 682
$ws.Cells.Item(17, 9).Value = 84
$ws.Cells.Item(18, 28).Value = 1016
$ws.Cells.Item(18, 9).Value = 114
$ws.Cells.Item(19, 28).Value = 1301
$ws.Cells.Item(19, 9).Value = 178
$ws.Cells.Item(20, 28).Value = 1582
$ws.Cells.Item(20, 9).Value = 281
$ws.Cells.Item(21, 28).Value = 1830
$ws.Cells.Item(21, 9).Value = 374
$ws.Cells.Item(22, 28).Value = 2388
$ws.Cells.Item(22, 9).Value = 461
$ws.Cells.Item(23, 28).Value = 3004
$ws.Cells.Item(23, 9).Value = 588
$ws.Cells.Item(24, 28).Value = 4141
$ws.Cells.Item(24, 9).Value = 786
$ws.Cells.Item(25, 28).Value = 5435
$ws.Cells.Item(25, 9).Value = 1026
$ws.Cells.Item(26, 28).Value = 6521
$ws.Cells.Item(26, 9).Value = 1152
$ws.Cells.Item(27, 28).Value = 7345
$ws.Cells.Item(27, 9).Value = 1281
$ws.Cells.Item(28, 28).Value = 7947
$ws.Cells.Item(28, 9).Value = 1424
$ws.Cells.Item(29, 28).Value = 9178
$ws.Cells.Item(29, 9).Value = 1571
$ws.Cells.Item(30, 28).Value = 10181
$ws.Cells.Item(30, 9).Value = 1673
$ws.Cells.Item(31, 28).Value = 11198
$ws.Cells.Item(31, 9).Value = 1796
$ws.Cells.Item(32, 28).Value = 12433
$ws.Cells.Item(32, 9).Value = 2021
$ws.Cells.Item(33, 28).Value = 13718
$ws.Cells.Item(33, 9).Value = 2209
$ws.Cells.Item(34, 28).Value = 14637
$ws.Cells.Item(34, 9).Value = 2405
$ws.Cells.Item(35, 28).Value = 15378
$ws.Cells.Item(35, 9).Value = 2511
$ws.Cells.Item(36, 28).Value = 16404
$ws.Cells.Item(36, 9).Value = 2718
$ws.Cells.Item(37, 28).Value = 17397
$ws.Cells.Item(37, 9).Value = 2948
$ws.Cells.Item(38, 28).Value = 18394
$ws.Cells.Item(38, 9).Value = 3105
$ws.Cells.Item(39, 28).Value = 19498
$ws.Cells.Item(39, 9).Value = 3316
$ws.Cells.Item(40, 28).Value = 20417
$ws.Cells.Item(40, 9).Value = 3503
$ws.Cells.Item(41, 28).Value = 21034
$ws.Cells.Item(41, 9).Value = 3630
$ws.Cells.Item(42, 28).Value = 21516
$ws.Cells.Item(42, 9).Value = 3695
$ws.Cells.Item(43, 28).Value = 22164
$ws.Cells.Item(43, 9).Value = 3814
$ws.Cells.Item(44, 28).Value = 22822
$ws.Cells.Item(44, 9).Value = 3974
$ws.Cells.Item(45, 28).Value = 23493
$ws.Cells.Item(45, 9).Value = 4105
$ws.Cells.Item(46, 28).Value = 24189
$ws.Cells.Item(46, 9).Value = 4198
$ws.Cells.Item(47, 28).Value = 24691
$ws.Cells.Item(47, 9).Value = 4240
$ws.Cells.Item(48, 28).Value = 25160
$ws.Cells.Item(48, 9).Value = 4289
$ws.Cells.Item(49, 28).Value = 25495
$ws.Cells.Item(49, 9).Value = 4311
$ws.Cells.Item(50, 28).Value = 25749
$ws.Cells.Item(50, 9).Value = 4361
$ws.Cells.Item(51, 28).Value = 26058
$ws.Cells.Item(51, 9).Value = 4414
$ws.Cells.Item(52, 28).Value = 26380
$ws.Cells.Item(52, 9).Value = 4480
$ws.Cells.Item(53, 28).Value = 26711
$ws.Cells.Item(53, 9).Value = 4545
$ws.Cells.Item(54, 28).Value = 27026
$ws.Cells.Item(54, 9).Value = 4605
$ws.Cells.Item(55, 28).Value = 27361
$ws.Cells.Item(55, 9).Value = 4654
$ws.Cells.Item(56, 28).Value = 27562
$ws.Cells.Item(56, 9).Value = 4665
$ws.Cells.Item(57, 28).Value = 27749
$ws.Cells.Item(57, 9).Value = 4694
$ws.Cells.Item(58, 28).Value = 27860
$ws.Cells.Item(58, 9).Value = 4710
$ws.Cells.Item(59, 28).Value = 27987

$ws = $wb.Worksheets.Item("Fatalities")
$ws.Cells.Item(29, 28).Value = 138
$ws.Cells.Item(29, 9).Value = 14
$ws.Cells.Item(30, 28).Value = 159
$ws.Cells.Item(30, 9).Value = 15
$ws.Cells.Item(31, 28).Value = 191
$ws.Cells.Item(31, 9).Value = 22
$ws.Cells.Item(32, 28).Value = 234
$ws.Cells.Item(32, 9).Value = 24
$ws.Cells.Item(33, 28).Value = 267
$ws.Cells.Item(33, 9).Value = 31
$ws.Cells.Item(34, 28).Value = 309
$ws.Cells.Item(34, 9).Value = 39
$ws.Cells.Item(35, 28).Value = 345
$ws.Cells.Item(35, 9).Value = 46
$ws.Cells.Item(36, 28).Value = 401
$ws.Cells.Item(36, 9).Value = 55
$ws.Cells.Item(37, 28).Value = 473
$ws.Cells.Item(37, 9).Value = 63
$ws.Cells.Item(38, 28).Value = 530
$ws.Cells.Item(38, 9).Value = 71
$ws.Cells.Item(39, 28).Value = 583
$ws.Cells.Item(39, 9).Value = 79
$ws.Cells.Item(40, 28).Value = 645
$ws.Cells.Item(40, 9).Value = 85
$ws.Cells.Item(41, 28).Value = 717
$ws.Cells.Item(41, 9).Value = 98
$ws.Cells.Item(42, 28).Value = 764
$ws.Cells.Item(42, 9).Value = 105
$ws.Cells.Item(43, 28).Value = 818
$ws.Cells.Item(43, 9).Value = 116
$ws.Cells.Item(44, 28).Value = 876
$ws.Cells.Item(44, 9).Value = 127
$ws.Cells.Item(45, 28).Value = 949
$ws.Cells.Item(45, 9).Value = 140
$ws.Cells.Item(46, 28).Value = 1002
$ws.Cells.Item(46, 9).Value = 149
$ws.Cells.Item(47, 28).Value = 1060
$ws.Cells.Item(47, 9).Value = 159
$ws.Cells.Item(48, 28).Value = 1099
$ws.Cells.Item(48, 9).Value = 167
$ws.Cells.Item(49, 28).Value = 1148
$ws.Cells.Item(49, 9).Value = 173
$ws.Cells.Item(50, 28).Value = 1175
$ws.Cells.Item(50, 9).Value = 176
$ws.Cells.Item(51, 28).Value = 1220
$ws.Cells.Item(51, 9).Value = 180
$ws.Cells.Item(52, 28).Value = 1270
$ws.Cells.Item(52, 9).Value = 184
$ws.Cells.Item(53, 28).Value = 1313
$ws.Cells.Item(53, 9).Value = 190
$ws.Cells.Item(54, 28).Value = 1359
$ws.Cells.Item(54, 9).Value = 196
$ws.Cells.Item(55, 28).Value = 1390
$ws.Cells.Item(55, 9).Value = 199
$ws.Cells.Item(56, 28).Value = 1415
$ws.Cells.Item(56, 9).Value = 201
$ws.Cells.Item(57, 28).Value = 1460
$ws.Cells.Item(57, 9).Value = 202
$ws.Cells.Item(58, 28).Value = 1490
$ws.Cells.Item(58, 9).Value = 205
$ws.Cells.Item(59, 28).Value = 1499

$ws = $wb.Worksheets.Item("Hospitalized")
$ws.Cells.Item(59, 28).Value = 1421
$ws.Cells.Item(59, 9).Value = 288

$ws = $wb.Worksheets.Item("ICU")
$ws.Cells.Item(59, 28).Value = 237
$ws.Cells.Item(59, 9).Value = 30

$ws = $wb.Worksheets.Item("Ventilated")
$ws.Cells.Item(59, 28).Value = 151
$ws.Cells.Item(59, 9).Value = 27

$ws = $wb.Worksheets.Item("Released")
$ws.Cells.Item(59, 28).Value = 4485
$ws.Cells.Item(59, 9).Value = 589

$ws = $wb.Worksheets.Item("Tested")
$ws.Cells.Item(5, 28).Value = 534
$ws.Cells.Item(5, 9).Value = 534
$ws.Cells.Item(6, 28).Value = 674
$ws.Cells.Item(6, 9).Value = 674
$ws.Cells.Item(7, 28).Value = 813
$ws.Cells.Item(7, 9).Value = 783
$ws.Cells.Item(8, 28).Value = 1136
$ws.Cells.Item(8, 9).Value = 871
$ws.Cells.Item(9, 28).Value = 1246
$ws.Cells.Item(9, 9).Value = 980
$ws.Cells.Item(10, 28).Value = 1355
$ws.Cells.Item(10, 9).Value = 1087
$ws.Cells.Item(11, 28).Value = 1421
$ws.Cells.Item(11, 9).Value = 1153
$ws.Cells.Item(12, 28).Value = 1522
$ws.Cells.Item(12, 9).Value = 1251
$ws.Cells.Item(14, 28).Value = 1669
$ws.Cells.Item(14, 9).Value = 1398
$ws.Cells.Item(15, 28).Value = 1782
$ws.Cells.Item(15, 9).Value = 1511
$ws.Cells.Item(16, 28).Value = 1992
$ws.Cells.Item(16, 9).Value = 1721
$ws.Cells.Item(17, 28).Value = 2243
$ws.Cells.Item(17, 9).Value = 1972
$ws.Cells.Item(18, 28).Value = 2616
$ws.Cells.Item(18, 9).Value = 2345
$ws.Cells.Item(19, 28).Value = 3098
$ws.Cells.Item(19, 9).Value = 2817
$ws.Cells.Item(20, 28).Value = 3742
$ws.Cells.Item(20, 9).Value = 3461
$ws.Cells.Item(21, 28).Value = 4177
$ws.Cells.Item(21, 9).Value = 3896
$ws.Cells.Item(22, 28).Value = 4874
$ws.Cells.Item(22, 9).Value = 4347
$ws.Cells.Item(23, 28).Value = 5506
$ws.Cells.Item(23, 9).Value = 4949
$ws.Cells.Item(24, 28).Value = 6299
$ws.Cells.Item(24, 9).Value = 5657
$ws.Cells.Item(25, 28).Value = 6894
$ws.Cells.Item(25, 9).Value = 6252
$ws.Cells.Item(26, 28).Value = 7301
$ws.Cells.Item(26, 9).Value = 6659
$ws.Cells.Item(27, 28).Value = 7669
$ws.Cells.Item(27, 9).Value = 7027
$ws.Cells.Item(28, 28).Value = 7910
$ws.Cells.Item(28, 9).Value = 7268
$ws.Cells.Item(29, 28).Value = 8306
$ws.Cells.Item(29, 9).Value = 7664
$ws.Cells.Item(30, 28).Value = 8654
$ws.Cells.Item(30, 9).Value = 8012
$ws.Cells.Item(31, 28).Value = 9087
$ws.Cells.Item(31, 9).Value = 8445
$ws.Cells.Item(32, 28).Value = 9600
$ws.Cells.Item(32, 9).Value = 8958
$ws.Cells.Item(33, 28).Value = 10092
$ws.Cells.Item(33, 9).Value = 9450
$ws.Cells.Item(34, 28).Value = 10699
$ws.Cells.Item(34, 9).Value = 10057
$ws.Cells.Item(35, 28).Value = 11008
$ws.Cells.Item(35, 9).Value = 10366
$ws.Cells.Item(36, 28).Value = 11668
$ws.Cells.Item(36, 9).Value = 11026
$ws.Cells.Item(37, 28).Value = 12426
$ws.Cells.Item(37, 9).Value = 11784
$ws.Cells.Item(38, 28).Value = 13091
$ws.Cells.Item(38, 9).Value = 12449
$ws.Cells.Item(39, 28).Value = 14043
$ws.Cells.Item(39, 9).Value = 13401
$ws.Cells.Item(40, 28).Value = 14881
$ws.Cells.Item(40, 9).Value = 14239
$ws.Cells.Item(41, 28).Value = 15475
$ws.Cells.Item(41, 9).Value = 14833
$ws.Cells.Item(42, 28).Value = 15779
$ws.Cells.Item(42, 9).Value = 15137
$ws.Cells.Item(43, 28).Value = 16434
$ws.Cells.Item(43, 9).Value = 15792
$ws.Cells.Item(44, 28).Value = 17285
$ws.Cells.Item(44, 9).Value = 16643
$ws.Cells.Item(45, 28).Value = 18106
$ws.Cells.Item(45, 9).Value = 17464
$ws.Cells.Item(46, 28).Value = 18628
$ws.Cells.Item(46, 9).Value = 17986
$ws.Cells.Item(47, 28).Value = 19048
$ws.Cells.Item(47, 9).Value = 18406
$ws.Cells.Item(48, 28).Value = 19417
$ws.Cells.Item(48, 9).Value = 18775
$ws.Cells.Item(49, 28).Value = 19596
$ws.Cells.Item(49, 9).Value = 18954
$ws.Cells.Item(50, 28).Value = 19907
$ws.Cells.Item(50, 9).Value = 19265
$ws.Cells.Item(51, 28).Value = 20462
$ws.Cells.Item(51, 9).Value = 19820
$ws.Cells.Item(52, 28).Value = 20938
$ws.Cells.Item(52, 9).Value = 20296
$ws.Cells.Item(53, 28).Value = 21630
$ws.Cells.Item(53, 9).Value = 20988
$ws.Cells.Item(54, 28).Value = 22100
$ws.Cells.Item(54, 9).Value = 21458
$ws.Cells.Item(55, 28).Value = 22552
$ws.Cells.Item(55, 9).Value = 21910
$ws.Cells.Item(56, 28).Value = 22743
$ws.Cells.Item(56, 9).Value = 22101
$ws.Cells.Item(57, 28).Value = 23135
$ws.Cells.Item(57, 9).Value = 22493
$ws.Cells.Item(58, 28).Value = 23377
$ws.Cells.Item(58, 9).Value = 22735
$ws.Cells.Item(59, 28).Value = 23377

Write-Output "Applied 283 cell updates across 7 sheets"